$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '34.493.03'
$ws.Range('E2').Value = '  +12.52%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.830.52'
$ws.Range('E3').Value = '  +9.20%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.21%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '230.65'
$ws.Range('E6').Value = '  +4.02%  '
$ws.Range('E7').Value = '  +0.22%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '31.70'
$ws.Range('E8').Value = '  +6.85%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '47.08'
$ws.Range('E9').Value = '  +6.65%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.285'
$ws.Range('E10').Value = '  +7.62%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0675'
$ws.Range('E11').Value = '  +3.88%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0934'
$ws.Range('E12').Value = '  +3.19%  '
$ws.Range('E13').Value = '  +9.22%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.833.97'
$ws.Range('E14').Value = '  +8.71%  '
$ws.Range('E15').Value = '  +6.14%  '
$ws.Range('B16').Value = 'Chainlink'
$ws.Range('C16').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '10.46'
$ws.Range('E16').Value = '  +1.81%  '
$ws.Range('B17').Value = 'WrappedBTC'
$ws.Range('C17').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '34.438.00'
$ws.Range('E17').Value = '  +12.33%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '4.27'
$ws.Range('E18').Value = '  +6.30%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '70.01'
$ws.Range('E19').Value = '  +5.43%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '260.85'
$ws.Range('E20').Value = '  +6.94%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0₃0754'
$ws.Range('E21').Value = '  +3.91%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.999'
$ws.Range('E22').Value = '  +0.16%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.58'
$ws.Range('E23').Value = '  +5.76%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '4.37'
$ws.Range('E24').Value = '  +2.37%  '
$ws.Range('E25').Value = '  +2.26%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '158.58'
$ws.Range('E26').Value = '  -0.07%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '16.76'
$ws.Range('E27').Value = '  +5.52%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.14'
$ws.Range('E28').Value = '  +6.66%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.116'
$ws.Range('E29').Value = '  +2.81%  '
$ws.Range('E30').Value = '  +0.24%  '
$ws.Range('E31').Value = '  +12.51%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0521'
$ws.Range('E32').Value = '  +5.15%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.22'
$ws.Range('E33').Value = '  +6.26%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.60'
$ws.Range('E34').Value = '  +9.36%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.552.08'
$ws.Range('E35').Value = '  +4.33%  '
$ws.Range('E36').Value = '  +1.95%  '
$ws.Range('E37').Value = '  +5.53%  '
$ws.Range('E38').Value = '  +217.47%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.647'
$ws.Range('E39').Value = '  +7.74%  '
$ws.Range('E40').Value = '  +6.52%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '85.17'
$ws.Range('E41').Value = '  +0.05%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.80'
$ws.Range('E42').Value = '  +4.86%  '
$ws.Range('E43').Value = '  +2.15%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.916'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.14'
$ws.Range('E45').Value = '  +8.27%  '
$ws.Range('E46').Value = '  +5.37%  '
$ws.Range('E47').Value = '  +5.95%  '
$ws.Range('E48').Value = '  +9.55%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '12.51'
$ws.Range('E49').Value = '  +28.21%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '5.84'
$ws.Range('E50').Value = '  +5.85%  '
$ws.Range('E51').Value = '  +0.26%  '
